$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 12, shifting rows 12-19 down to 13-20.
$ws.Rows.Item(12).Insert()

# New row 12 content: label + ratio value.
$ws.Range("A12").Value = "Ratio of too-sick-to-vote to contagious voters"
$ws.Range("B12").Value = 0.5

# Update the "Expected voters at polling place" formula (now row 18) to include the new factor.
$ws.Range("B18").Formula = "=B7*(1-B14)*B13*B12"

$ws.Range("A5").Select()
